# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @{
    2  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    3  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    4  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    5  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    6  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    7  = @{ B = 0.7287194209349384; C = 1.65323645889881;  D = 3.082599426703578;  E = 0.4998867070740569; G = 5.964442013611383 }
    8  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    9  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 3.082599426703578;  E = 0.4998867070740569; G = 8.418600821238126 }
    10 = @{ B = 1.505614041169197;  C = 9.226618575922256; D = 0.7127328510149897; E = 6.48142807727062;  G = 17.92639354537706 }
    11 = @{ B = 0.7287194209349384; C = 1.65323645889881;  D = 0.7127328510149897; E = 6.48142807727062;  G = 9.576116808119359 }
    12 = @{ B = 0.7287194209349384; C = 1766.335244827366; D = 157.8057217802531;  E = 6.48142807727062;  G = 1931.351114105825 }
    13 = @{ B = 0.1554434735375247; C = 0.3375848360084654; D = 0.1529057820181812; E = 6.48142807727062; G = 7.127362168834791 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 2).Value = $vals.B
    $ws.Cells.Item($r, 3).Value = $vals.C
    $ws.Cells.Item($r, 4).Value = $vals.D
    $ws.Cells.Item($r, 5).Value = $vals.E
    $ws.Cells.Item($r, 7).Value = $vals.G
}
